# Adds the new "checkViewProjects" sheet (View Projects in Dashboard test),
# mirroring the pattern used by the existing "checkArtifactsManager" sheet.

$wb = $excel.ActiveWorkbook

# The new sheet is appended right after the current last sheet
# (checkArtifactsManager), becoming the new active/selected tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "checkViewProjects"

# Column A holds the title key/value pair, same shape as the other
# "checkXxx" sheets (e.g. checkArtifactsManager -> A1/A2).
$newSheet.Range("A1").Value = "viewProjectsTitle"
$newSheet.Range("A2").Value = "View Projects"

# Match the column width used on these small "title" sheets.
$newSheet.Columns.Item(1).ColumnWidth = 28

# Leave selection on A3, as captured in the saved workbook state.
$newSheet.Range("A3").Select()
